$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.473.18'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.909.01'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.87'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4838'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4058'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08167'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.39'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.46%  '
$ws.Range('D12').Value = '1.939.18'
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.011'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.161'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.38'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '29.488.66'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.627'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.75'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.191'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').Value = '2.178.13'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.94'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.462'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +7.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.112'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.34'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.029'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09519'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.513'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.562'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.390'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02276'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06111'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.176'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.85'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5951'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.977'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.279'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.369'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.91%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.49'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.67%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07611'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5568'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.948'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.64'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('E51').Value = '  +2.35%  '
